# Auto-generated edit script applying cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "28.423.48"
$ws.Cells.Item(2, 5).Value = "  -0.23%  "
$ws.Cells.Item(3, 4).Value = "1.573.24"
$ws.Cells.Item(3, 5).Value = "  -0.03%  "
$ws.Cells.Item(4, 5).Value = "  -0.30%  "
$ws.Cells.Item(5, 4).Value = "'211.98"
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).Value = "  -0.16%  "
$ws.Cells.Item(6, 4).Value = "'0.491"
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).Value = "  -0.17%  "
$ws.Cells.Item(7, 5).Value = "  -0.30%  "
$ws.Cells.Item(8, 4).Value = "'44.56"
$ws.Cells.Item(8, 4).ClearFormats()
$ws.Cells.Item(8, 5).Value = "  -5.25%  "
$ws.Cells.Item(9, 4).Value = "'23.72"
$ws.Cells.Item(9, 4).ClearFormats()
$ws.Cells.Item(9, 5).Value = "  -1.94%  "
$ws.Cells.Item(10, 5).Value = "  -0.71%  "
$ws.Cells.Item(11, 5).Value = "  -0.72%  "
$ws.Cells.Item(12, 4).Value = "'0.0893"
$ws.Cells.Item(12, 4).ClearFormats()
$ws.Cells.Item(12, 5).Value = "  +1.15%  "
$ws.Cells.Item(13, 4).Value = "1.798.90"
$ws.Cells.Item(13, 5).Value = "  -0.01%  "
$ws.Cells.Item(14, 4).Value = "1.571.55"
$ws.Cells.Item(14, 5).Value = "  -0.36%  "
$ws.Cells.Item(15, 5).Value = "  -0.33%  "
$ws.Cells.Item(16, 4).Value = "28.406.17"
$ws.Cells.Item(16, 5).Value = "  -0.42%  "
$ws.Cells.Item(17, 5).Value = "  -1.20%  "
$ws.Cells.Item(18, 4).Value = "'61.64"
$ws.Cells.Item(18, 4).ClearFormats()
$ws.Cells.Item(18, 5).Value = "  -0.96%  "
$ws.Cells.Item(19, 4).Value = "'229.92"
$ws.Cells.Item(19, 4).ClearFormats()
$ws.Cells.Item(19, 5).Value = "  +0.98%  "
$ws.Cells.Item(20, 4).Value = "'7.40"
$ws.Cells.Item(20, 4).ClearFormats()
$ws.Cells.Item(20, 5).Value = "  +0.25%  "
$ws.Cells.Item(21, 4).Value = "0.0₃0684"
$ws.Cells.Item(21, 5).Value = "  -1.37%  "
$ws.Cells.Item(22, 5).Value = "  -0.13%  "
$ws.Cells.Item(23, 5).Value = "  +1.51%  "
$ws.Cells.Item(24, 4).Value = "'9.01"
$ws.Cells.Item(24, 4).ClearFormats()
$ws.Cells.Item(24, 5).Value = "  -1.30%  "
$ws.Cells.Item(25, 4).Value = "'2.04"
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Cells.Item(25, 5).Value = "  +1.29%  "
$ws.Cells.Item(26, 4).Value = "'151.35"
$ws.Cells.Item(26, 4).ClearFormats()
$ws.Cells.Item(26, 5).Value = "  -0.05%  "
$ws.Cells.Item(27, 4).Value = "'14.90"
$ws.Cells.Item(27, 4).ClearFormats()
$ws.Cells.Item(27, 5).Value = "  -0.55%  "
$ws.Cells.Item(28, 2).Value = "Stellar"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(28, 4).Value = "'0.104"
$ws.Cells.Item(28, 4).ClearFormats()
$ws.Cells.Item(28, 5).Value = "  -0.60%  "
$ws.Cells.Item(29, 2).Value = "Cosmos"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(29, 4).Value = "'6.36"
$ws.Cells.Item(29, 4).ClearFormats()
$ws.Cells.Item(29, 5).Value = "  -1.36%  "
$ws.Cells.Item(30, 5).Value = "  -0.28%  "
$ws.Cells.Item(31, 5).Value = "  +4.14%  "
$ws.Cells.Item(32, 5).Value = "  -3.37%  "
$ws.Cells.Item(33, 5).Value = "  -0.73%  "
$ws.Cells.Item(34, 4).Value = "'3.12"
$ws.Cells.Item(34, 4).ClearFormats()
$ws.Cells.Item(34, 5).Value = "  -0.88%  "
$ws.Cells.Item(35, 4).Value = "1.382.03"
$ws.Cells.Item(35, 5).Value = "  -0.78%  "
$ws.Cells.Item(36, 5).Value = "  +4.64%  "
$ws.Cells.Item(37, 4).Value = "'1.52"
$ws.Cells.Item(37, 4).ClearFormats()
$ws.Cells.Item(37, 5).Value = "  -1.96%  "
$ws.Cells.Item(38, 5).Value = "  -0.08%  "
$ws.Cells.Item(39, 5).Value = "  +0.22%  "
$ws.Cells.Item(40, 5).Value = "  -1.66%  "
$ws.Cells.Item(41, 5).Value = "  -2.47%  "
$ws.Cells.Item(42, 2).Value = "RenderToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(42, 4).Value = "'1.91"
$ws.Cells.Item(42, 4).ClearFormats()
$ws.Cells.Item(42, 5).Value = "  +2.84%  "
$ws.Cells.Item(43, 2).Value = "PaxDollar"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Cells.Item(43, 4).Value = "'1.00"
$ws.Cells.Item(43, 4).ClearFormats()
$ws.Cells.Item(43, 5).Value = "  -0.20%  "
$ws.Cells.Item(44, 5).Value = "  -0.65%  "
$ws.Cells.Item(45, 4).Value = "'0.0475"
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Cells.Item(45, 5).Value = "  +1.35%  "
$ws.Cells.Item(46, 4).Value = "'5.37"
$ws.Cells.Item(46, 4).ClearFormats()
$ws.Cells.Item(46, 5).Value = "  -4.45%  "
$ws.Cells.Item(47, 4).Value = "'62.38"
$ws.Cells.Item(47, 4).ClearFormats()
$ws.Cells.Item(47, 5).Value = "  -1.00%  "
$ws.Cells.Item(48, 2).Value = "RocketPoolETH"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Cells.Item(48, 4).Value = "1.710.70"
$ws.Cells.Item(48, 5).Value = "  -0.08%  "
$ws.Cells.Item(49, 2).Value = "WEMIXToken"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(49, 4).Value = "'0.919"
$ws.Cells.Item(49, 4).ClearFormats()
$ws.Cells.Item(49, 5).Value = "  -6.24%  "
$ws.Cells.Item(50, 4).Value = "'2.13"
$ws.Cells.Item(50, 4).ClearFormats()
$ws.Cells.Item(50, 5).Value = "  -0.35%  "
$ws.Cells.Item(51, 4).Value = "'85.27"
$ws.Cells.Item(51, 4).ClearFormats()
$ws.Cells.Item(51, 5).Value = "  -0.92%  "
